$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.316.08'
$ws.Range("E2").Value = '  -0.79%  '
$ws.Range("D3").Value = '1.638.28'
$ws.Range("E3").Value = '  -1.68%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.54'
$ws.Range("E5").Value = '  -1.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.536'
$ws.Range("E6").Value = '  +4.33%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.85'
$ws.Range("E8").Value = '  -3.39%  '
$ws.Range("E9").Value = '  -2.43%  '
$ws.Range("E10").Value = '  -2.26%  '
$ws.Range("E11").Value = '  +1.40%  '
$ws.Range("D12").Value = '1.870.73'
$ws.Range("E12").Value = '  -1.65%  '
$ws.Range("D13").Value = '1.632.75'
$ws.Range("E13").Value = '  -2.33%  '
$ws.Range("E14").Value = '  -3.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.556'
$ws.Range("E15").Value = '  +0.30%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.18'
$ws.Range("E16").Value = '  -3.26%  '
$ws.Range("D17").Value = '27.297.12'
$ws.Range("E17").Value = '  -0.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '226.56'
$ws.Range("E18").Value = '  -8.71%  '
$ws.Range("D19").Value = '0.0₃0716'
$ws.Range("E19").Value = '  -2.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.38'
$ws.Range("E20").Value = '  -2.51%  '
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("E22").Value = '  -4.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.24'
$ws.Range("E23").Value = '  -0.77%  '
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.13'
$ws.Range("E25").Value = '  +0.41%  '
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.114'
$ws.Range("E26").Value = '  +1.66%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.89'
$ws.Range("E27").Value = '  -3.65%  '
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.46'
$ws.Range("E29").Value = '  -6.99%  '
$ws.Range("E30").Value = '  -4.91%  '
$ws.Range("E31").Value = '  -5.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.25'
$ws.Range("E32").Value = '  -3.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.08'
$ws.Range("E33").Value = '  -1.22%  '
$ws.Range("D34").Value = '1.388.22'
$ws.Range("E34").Value = '  -5.58%  '
$ws.Range("E35").Value = '  -1.30%  '
$ws.Range("E36").Value = '  -0.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.875'
$ws.Range("E37").Value = '  -7.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.557'
$ws.Range("E38").Value = '  -3.47%  '
$ws.Range("E39").Value = '  -3.63%  '
$ws.Range("E40").Value = '  +0.73%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.47'
$ws.Range("E42").Value = '  -1.63%  '
$ws.Range("E43").Value = '  -0.39%  '
$ws.Range("E44").Value = '  -0.33%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.786'
$ws.Range("E45").Value = '  -0.69%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.91'
$ws.Range("E46").Value = '  -8.58%  '
$ws.Range("D47").Value = '1.781.22'
$ws.Range("E47").Value = '  -1.54%  '
$ws.Range("E48").Value = '  -3.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.86'
$ws.Range("E49").Value = '  -3.13%  '
$ws.Range("E50").Value = '  -4.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0974'
$ws.Range("E51").Value = '  -4.77%  '
